# This script reproduces the crypto price/volume refresh captured in the
# commit "Updated cryptos list ... with GitHub Actions": most rows keep the
# same coin/link, but their Price (column D) and Volume(1h) (column E) are
# refreshed, and the PEPE / Binance-PegBSC-USD rows (29-30) swap places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price strings such as "65.013.71" or "1.00" that look
# numeric but must stay plain text (to preserve formatting such as trailing
# zeros and the "thousand.thousand.decimal" style grouping). Temporarily force
# the whole data range in column D to Text format before writing, then restore
# the original (default/no explicit) cell style once all values are written so
# the saved file does not end up with stray numeric formatting applied.
$priceRange = $ws.Range("D2:D51")
$originalPriceStyle = $priceRange.Style
$priceRange.NumberFormat = "@"

$ws.Range("D2").Value = '65.013.71'
$ws.Range("E2").Value = '  +2.55%  '
$ws.Range("D3").Value = '2.635.70'
$ws.Range("E3").Value = '  +2.36%  '
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").Value = '596.99'
$ws.Range("E5").Value = '  +1.64%  '
$ws.Range("D6").Value = '155.77'
$ws.Range("E6").Value = '  +4.51%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '0.591'
$ws.Range("E8").Value = '  +1.07%  '
$ws.Range("E9").Value = '  +7.31%  '
$ws.Range("D10").Value = '0.401'
$ws.Range("E10").Value = '  +4.99%  '
$ws.Range("E11").Value = '  +1.05%  '
$ws.Range("D12").Value = '0.154'
$ws.Range("E12").Value = '  +2.37%  '
$ws.Range("E13").Value = '  +6.41%  '
$ws.Range("D14").Value = '0.0000186'
$ws.Range("E14").Value = '  +21.26%  '
$ws.Range("D15").Value = '3.112.02'
$ws.Range("E15").Value = '  +2.48%  '
$ws.Range("D16").Value = '64.946.20'
$ws.Range("E16").Value = '  +2.70%  '
$ws.Range("D17").Value = '2.640.20'
$ws.Range("E17").Value = '  +1.63%  '
$ws.Range("D18").Value = '12.55'
$ws.Range("E18").Value = '  +3.06%  '
$ws.Range("E19").Value = '  +2.60%  '
$ws.Range("D20").Value = '352.41'
$ws.Range("E20").Value = '  +2.09%  '
$ws.Range("D21").Value = '7.33'
$ws.Range("E21").Value = '  +7.53%  '
$ws.Range("E22").Value = '  +0.19%  '
$ws.Range("D23").Value = '68.24'
$ws.Range("E23").Value = '  +2.10%  '
$ws.Range("D24").Value = '1.69'
$ws.Range("E24").Value = '  +0.90%  '
$ws.Range("D25").Value = '9.53'
$ws.Range("E25").Value = '  +4.91%  '
$ws.Range("E26").Value = '  -1.36%  '
$ws.Range("D27").Value = '0.164'
$ws.Range("E27").Value = '  +1.04%  '
$ws.Range("D28").Value = '8.08'
$ws.Range("E28").Value = '  +0.95%  '
$ws.Range("B29").Value = 'Binance-PegBSC-USD'
$ws.Range("C29").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.00%  '
$ws.Range("B30").Value = 'PEPE'
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D30").Value = '0.0₃0947'
$ws.Range("E30").Value = '  +10.54%  '
$ws.Range("E31").Value = '  +4.06%  '
$ws.Range("D32").Value = '508.03'
$ws.Range("E32").Value = '  -8.30%  '
$ws.Range("D33").Value = '1.76'
$ws.Range("E33").Value = '  +1.04%  '
$ws.Range("E34").Value = '  +7.97%  '
$ws.Range("D35").Value = '6.34'
$ws.Range("E35").Value = '  +5.86%  '
$ws.Range("E36").Value = '  +2.89%  '
$ws.Range("D37").Value = '20.28'
$ws.Range("E37").Value = '  +4.59%  '
$ws.Range("D38").Value = '163.29'
$ws.Range("E38").Value = '  -1.18%  '
$ws.Range("E39").Value = '  +5.59%  '
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  +0.01%  '
$ws.Range("E41").Value = '  -0.02%  '
$ws.Range("D42").Value = '42.28'
$ws.Range("E42").Value = '  +6.44%  '
$ws.Range("D43").Value = '165.50'
$ws.Range("E43").Value = '  +0.27%  '
$ws.Range("D44").Value = '4.08'
$ws.Range("E44").Value = '  +2.88%  '
$ws.Range("E45").Value = '  +4.50%  '
$ws.Range("D46").Value = '23.08'
$ws.Range("E46").Value = '  +1.49%  '
$ws.Range("D47").Value = '2.20'
$ws.Range("E47").Value = '  +7.87%  '
$ws.Range("D48").Value = '0.646'
$ws.Range("E48").Value = '  +3.28%  '
$ws.Range("D49").Value = '0.0255'
$ws.Range("E49").Value = '  +2.43%  '
$ws.Range("D50").Value = '0.0983'
$ws.Range("E50").Value = '  +2.50%  '
$ws.Range("D51").Value = '19.48'
$ws.Range("E51").Value = '  +2.67%  '

$priceRange.Style = $originalPriceStyle
